$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New summary rows ---------------------------------------------------

# Row 12: average of the |S*|/n column (J)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Rows 14-17: labelled summary statistics in columns A (label) / B (value)
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Taller rows for 14-17 (matches ht="15.6" in the saved file)
$ws.Range("A14:B17").RowHeight = 15.6

# --- Formatting for the new summary values (bold, 12pt, vertically centred)
# Build the style once on a scratch cell, then fan it out with a format-only
# paste so we don't leave a trail of intermediate cell styles behind.
$scratch = $ws.Range("ZZ1")
$scratch.Font.Bold = $true
$scratch.Font.Size = 12
$scratch.VerticalAlignment = -4108

$scratch.Copy() | Out-Null
$ws.Range("B14:B17").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false
$scratch.Clear() | Out-Null

# --- Page setup -----------------------------------------------------------
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# --- Selection / view state -------------------------------------------
$ws.Range("A14:B17").Select()
